# Update "想去人数" (want-to-go count) values that were refreshed by the
# gh-pages data generation bot (commit: "Update gh-pages to output
# generated at 456a3b4").
#
# Sheet "展览" (Exhibitions, sheet1):
#   F2: 1050 -> 1049
#   F3: 366  -> 373
#   F4: 2987 -> 3007
#   F6: 624  -> 626
#
# Sheet "全部类型" (All Types, sheet4) mirrors the same events two rows
# further down:
#   F4: 1050 -> 1049
#   F5: 366  -> 373
#   F6: 2987 -> 3007
#   F8: 624  -> 626

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 1049
$wsExhibition.Range("F3").Value = 373
$wsExhibition.Range("F4").Value = 3007
$wsExhibition.Range("F6").Value = 626

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1049
$wsAll.Range("F5").Value = 373
$wsAll.Range("F6").Value = 3007
$wsAll.Range("F8").Value = 626
